# Natmi following Dr Hou advice
# Rewrite LR-pair table rows 2-10 (Sema3a -> Plxna1, senders ECs/FAPs/sCs x targets ECs/FAPs/sCs)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> ECs
$ws.Cells.Item(2,1).Value = "ECs"
$ws.Cells.Item(2,2).Value = "Sema3a"
$ws.Cells.Item(2,3).Value = "Plxna1"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 2
$ws.Cells.Item(2,6).Value = 0.6666666666666666
$ws.Cells.Item(2,7).Value = 0.5674196666666667
$ws.Cells.Item(2,8).Value = 1.702259
$ws.Cells.Item(2,9).Value = 0.07864125446886469
$ws.Cells.Item(2,10).Value = 0.07864125446886468
$ws.Cells.Item(2,11).Value = 3
$ws.Cells.Item(2,12).Value = 1
$ws.Cells.Item(2,13).Value = 7.023694333333334
$ws.Cells.Item(2,14).Value = 21.071083
$ws.Cells.Item(2,15).Value = 0.1590811435055747
$ws.Cells.Item(2,16).Value = 0.1590811435055747
$ws.Cells.Item(2,17).Value = 3.985382297388556
$ws.Cells.Item(2,18).Value = 35.868440676497
$ws.Cells.Item(2,19).Value = 0.01251034068761988
$ws.Cells.Item(2,20).Value = 0.01251034068761988

# Row 3: ECs -> FAPs
$ws.Cells.Item(3,1).Value = "ECs"
$ws.Cells.Item(3,2).Value = "Sema3a"
$ws.Cells.Item(3,3).Value = "Plxna1"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = 2
$ws.Cells.Item(3,6).Value = 0.6666666666666666
$ws.Cells.Item(3,7).Value = 0.5674196666666667
$ws.Cells.Item(3,8).Value = 1.702259
$ws.Cells.Item(3,9).Value = 0.07864125446886469
$ws.Cells.Item(3,10).Value = 0.07864125446886468
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 13.91445766666667
$ws.Cells.Item(3,14).Value = 41.74337300000001
$ws.Cells.Item(3,15).Value = 0.3151515045818828
$ws.Cells.Item(3,16).Value = 0.3151515045818827
$ws.Cells.Item(3,17).Value = 7.895336931067446
$ws.Cells.Item(3,18).Value = 71.05803237960701
$ws.Cells.Item(3,19).Value = 0.02478390966806943
$ws.Cells.Item(3,20).Value = 0.02478390966806941

# Row 4: ECs -> sCs
$ws.Cells.Item(4,1).Value = "ECs"
$ws.Cells.Item(4,2).Value = "Sema3a"
$ws.Cells.Item(4,3).Value = "Plxna1"
$ws.Cells.Item(4,4).Value = "sCs"
$ws.Cells.Item(4,5).Value = 2
$ws.Cells.Item(4,6).Value = 0.6666666666666666
$ws.Cells.Item(4,7).Value = 0.5674196666666667
$ws.Cells.Item(4,8).Value = 1.702259
$ws.Cells.Item(4,9).Value = 0.07864125446886469
$ws.Cells.Item(4,10).Value = 0.07864125446886468
$ws.Cells.Item(4,11).Value = 3
$ws.Cells.Item(4,12).Value = 1
$ws.Cells.Item(4,13).Value = 23.213494
$ws.Cells.Item(4,14).Value = 69.640482
$ws.Cells.Item(4,15).Value = 0.5257673519125425
$ws.Cells.Item(4,16).Value = 0.5257673519125424
$ws.Cells.Item(4,17).Value = 13.17179302764867
$ws.Cells.Item(4,18).Value = 118.546137248838
$ws.Cells.Item(4,19).Value = 0.04134700411317539
$ws.Cells.Item(4,20).Value = 0.04134700411317537

# Row 5: FAPs -> ECs
$ws.Cells.Item(5,1).Value = "FAPs"
$ws.Cells.Item(5,2).Value = "Sema3a"
$ws.Cells.Item(5,3).Value = "Plxna1"
$ws.Cells.Item(5,4).Value = "ECs"
$ws.Cells.Item(5,5).Value = 3
$ws.Cells.Item(5,6).Value = 1
$ws.Cells.Item(5,7).Value = 0.7227763333333334
$ws.Cells.Item(5,8).Value = 2.168329
$ws.Cells.Item(5,9).Value = 0.1001728366019618
$ws.Cells.Item(5,10).Value = 0.1001728366019618
$ws.Cells.Item(5,11).Value = 3
$ws.Cells.Item(5,12).Value = 1
$ws.Cells.Item(5,13).Value = 7.023694333333334
$ws.Cells.Item(5,14).Value = 21.071083
$ws.Cells.Item(5,15).Value = 0.1590811435055747
$ws.Cells.Item(5,16).Value = 0.1590811435055747
$ws.Cells.Item(5,17).Value = 5.076560036700778
$ws.Cells.Item(5,18).Value = 45.689040330307
$ws.Cells.Item(5,19).Value = 0.01593560939483717
$ws.Cells.Item(5,20).Value = 0.01593560939483717

# Row 6: FAPs -> FAPs
$ws.Cells.Item(6,1).Value = "FAPs"
$ws.Cells.Item(6,2).Value = "Sema3a"
$ws.Cells.Item(6,3).Value = "Plxna1"
$ws.Cells.Item(6,4).Value = "FAPs"
$ws.Cells.Item(6,5).Value = 3
$ws.Cells.Item(6,6).Value = 1
$ws.Cells.Item(6,7).Value = 0.7227763333333334
$ws.Cells.Item(6,8).Value = 2.168329
$ws.Cells.Item(6,9).Value = 0.1001728366019618
$ws.Cells.Item(6,10).Value = 0.1001728366019618
$ws.Cells.Item(6,11).Value = 3
$ws.Cells.Item(6,12).Value = 1
$ws.Cells.Item(6,13).Value = 13.91445766666667
$ws.Cells.Item(6,14).Value = 41.74337300000001
$ws.Cells.Item(6,15).Value = 0.3151515045818828
$ws.Cells.Item(6,16).Value = 0.3151515045818827
$ws.Cells.Item(6,17).Value = 10.05704069263522
$ws.Cells.Item(6,18).Value = 90.51336623371701
$ws.Cells.Item(6,19).Value = 0.03156962017334337
$ws.Cells.Item(6,20).Value = 0.03156962017334335

# Row 7: FAPs -> sCs
$ws.Cells.Item(7,1).Value = "FAPs"
$ws.Cells.Item(7,2).Value = "Sema3a"
$ws.Cells.Item(7,3).Value = "Plxna1"
$ws.Cells.Item(7,4).Value = "sCs"
$ws.Cells.Item(7,5).Value = 3
$ws.Cells.Item(7,6).Value = 1
$ws.Cells.Item(7,7).Value = 0.7227763333333334
$ws.Cells.Item(7,8).Value = 2.168329
$ws.Cells.Item(7,9).Value = 0.1001728366019618
$ws.Cells.Item(7,10).Value = 0.1001728366019618
$ws.Cells.Item(7,11).Value = 3
$ws.Cells.Item(7,12).Value = 1
$ws.Cells.Item(7,13).Value = 23.213494
$ws.Cells.Item(7,14).Value = 69.640482
$ws.Cells.Item(7,15).Value = 0.5257673519125425
$ws.Cells.Item(7,16).Value = 0.5257673519125424
$ws.Cells.Item(7,17).Value = 16.77816407717533
$ws.Cells.Item(7,18).Value = 151.003476694578
$ws.Cells.Item(7,19).Value = 0.05266760703378127
$ws.Cells.Item(7,20).Value = 0.05266760703378126

# Row 8: sCs -> ECs
$ws.Cells.Item(8,1).Value = "sCs"
$ws.Cells.Item(8,2).Value = "Sema3a"
$ws.Cells.Item(8,3).Value = "Plxna1"
$ws.Cells.Item(8,4).Value = "ECs"
$ws.Cells.Item(8,5).Value = 3
$ws.Cells.Item(8,6).Value = 1
$ws.Cells.Item(8,7).Value = 5.925096666666668
$ws.Cells.Item(8,8).Value = 17.77529
$ws.Cells.Item(8,9).Value = 0.8211859089291735
$ws.Cells.Item(8,10).Value = 0.8211859089291734
$ws.Cells.Item(8,11).Value = 3
$ws.Cells.Item(8,12).Value = 1
$ws.Cells.Item(8,13).Value = 7.023694333333334
$ws.Cells.Item(8,14).Value = 21.071083
$ws.Cells.Item(8,15).Value = 0.1590811435055747
$ws.Cells.Item(8,16).Value = 0.1590811435055747
$ws.Cells.Item(8,17).Value = 41.6160678821189
$ws.Cells.Item(8,18).Value = 374.5446109390701
$ws.Cells.Item(8,19).Value = 0.1306351934231177
$ws.Cells.Item(8,20).Value = 0.1306351934231176

# Row 9: sCs -> FAPs
$ws.Cells.Item(9,1).Value = "sCs"
$ws.Cells.Item(9,2).Value = "Sema3a"
$ws.Cells.Item(9,3).Value = "Plxna1"
$ws.Cells.Item(9,4).Value = "FAPs"
$ws.Cells.Item(9,5).Value = 3
$ws.Cells.Item(9,6).Value = 1
$ws.Cells.Item(9,7).Value = 5.925096666666668
$ws.Cells.Item(9,8).Value = 17.77529
$ws.Cells.Item(9,9).Value = 0.8211859089291735
$ws.Cells.Item(9,10).Value = 0.8211859089291734
$ws.Cells.Item(9,11).Value = 3
$ws.Cells.Item(9,12).Value = 1
$ws.Cells.Item(9,13).Value = 13.91445766666667
$ws.Cells.Item(9,14).Value = 41.74337300000001
$ws.Cells.Item(9,15).Value = 0.3151515045818828
$ws.Cells.Item(9,16).Value = 0.3151515045818827
$ws.Cells.Item(9,17).Value = 82.44450673924113
$ws.Cells.Item(9,18).Value = 742.0005606531702
$ws.Cells.Item(9,19).Value = 0.2587979747404701
$ws.Cells.Item(9,20).Value = 0.2587979747404699

# Row 10: sCs -> sCs
$ws.Cells.Item(10,1).Value = "sCs"
$ws.Cells.Item(10,2).Value = "Sema3a"
$ws.Cells.Item(10,3).Value = "Plxna1"
$ws.Cells.Item(10,4).Value = "sCs"
$ws.Cells.Item(10,5).Value = 3
$ws.Cells.Item(10,6).Value = 1
$ws.Cells.Item(10,7).Value = 5.925096666666668
$ws.Cells.Item(10,8).Value = 17.77529
$ws.Cells.Item(10,9).Value = 0.8211859089291735
$ws.Cells.Item(10,10).Value = 0.8211859089291734
$ws.Cells.Item(10,11).Value = 3
$ws.Cells.Item(10,12).Value = 1
$ws.Cells.Item(10,13).Value = 23.213494
$ws.Cells.Item(10,14).Value = 69.640482
$ws.Cells.Item(10,15).Value = 0.5257673519125425
$ws.Cells.Item(10,16).Value = 0.5257673519125424
$ws.Cells.Item(10,17).Value = 137.5421959210867
$ws.Cells.Item(10,18).Value = 1237.87976328978
$ws.Cells.Item(10,19).Value = 0.4317527407655859
$ws.Cells.Item(10,20).Value = 0.4317527407655857

